# ---------------------------------------------------------------------------
# Applies the "last update (safe backup)" commit to Overview_fullweek_DRres.xlsx
#
# Summary of the change:
#   - reserves: numeric values (cols C:N, rows 2-6) updated to new model run.
#   - capacities: unchanged.
#   - old "curtailment" sheet (PV/WIND_OFF/WIND_ON x curt/curt30/curt50/curt70)
#     is removed.
#   - old "storage" sheet (BEL_Z x STOR_L/STOR_M_C/STOR_M_E/STOR_S x
#     BEL/BEL30/BEL50/BEL70) is kept, gains a trailing STOR_S row, and is
#     moved right after "capacities".
#   - five new sheets are appended after it: generation, demandshift,
#     curtailment (new layout), storc, stordisc, maxshift.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: stamp the header style (bold font, thin border, centered/top
# aligned) used throughout this workbook (cellXfs index 1) onto a range by
# copying it off an already-styled cell -- this reuses the existing style
# instead of minting new ones in styles.xml.
# ---------------------------------------------------------------------------
$styleDonor = $wb.Worksheets.Item("reserves").Range("A1")
function Stamp-HeaderStyle($range) {
    $styleDonor.Copy()
    $range.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) reserves: refresh the numeric body (headers / labels untouched)
# ---------------------------------------------------------------------------
$reserves = $wb.Worksheets.Item("reserves")

$reservesData = @{
    2 = @("0.0013018120826643","0","1552040.93371274","696.835331348669","0","2638748.0651016","337046.459562756","2273917.40538827","1220670.65150236","68515.4357844881","50759.7164501673","4117407.65415917")
    3 = @("19926.8648783436","41517.7031028084","1516304.836956","149867.513627566","88018.2314367009","2427228.39482808","217443.35960642","2251050.94326732","1391082.82475575","16693.596305258","919152.888146678","3434260.45613125")
    4 = @("3860.80401347566","2221.93922025343","869919.656766282","48434.7047055267","3288.96132494043","824278.733969547","123233.591980567","32420.9857077419","720347.822311693","44936.515822082","101998.685544654","729067.19863322")
    5 = @("34.330083114472","0","9366743.59523794","13222.9232904404","0","12953259.1110675","7894874.10317724","1309734.19496392","7693855.07612112","19533.1781327815","12110226.6059994","6124162.82984752")
    6 = @("574472.740094029","274362.802142533","8761177.98945451","1931604.65251127","419372.917992239","10910710.7278905","8407972.81066769","5164480.15955406","3673238.89879365","3642350.09775439","15197443.0106239","688026.467582071")
}

foreach ($row in $reservesData.Keys) {
    $vals = $reservesData[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 3 + $i   # column C is index 3
        $reserves.Cells.Item($row, $col).Value = [double]$vals[$i]
    }
}

# ---------------------------------------------------------------------------
# 2) Rename the two sheets in place (preserves sheetId / position semantics)
#    curtailment(id=3) -> storage         (cleared, repopulated below)
#    storage(id=4)     -> generation      (cleared, repopulated below)
# ---------------------------------------------------------------------------
$futureStorage = $wb.Worksheets.Item("curtailment")
$futureGeneration = $wb.Worksheets.Item("storage")

$futureStorage.Cells.Clear() | Out-Null
$futureGeneration.Cells.Clear() | Out-Null

$futureStorage.Name = "storage_tmp_rename"
$futureGeneration.Name = "generation"
$futureStorage.Name = "storage"

# ---------------------------------------------------------------------------
# 3) Append the new sheets, in the exact target order, right after
#    "generation" (which is right after "storage").
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("generation")

$demandshift = $wb.Worksheets.Add($null, $anchor)
$demandshift.Name = "demandshift"

$curtailment = $wb.Worksheets.Add($null, $demandshift)
$curtailment.Name = "curtailment"

$storc = $wb.Worksheets.Add($null, $curtailment)
$storc.Name = "storc"

$stordisc = $wb.Worksheets.Add($null, $storc)
$stordisc.Name = "stordisc"

$maxshift = $wb.Worksheets.Add($null, $stordisc)
$maxshift.Name = "maxshift"

# ---------------------------------------------------------------------------
# 4) storage (renamed from curtailment) -- same layout as the old "storage"
#    sheet, with string indices shifted and a new trailing STOR_S row.
# ---------------------------------------------------------------------------
$storage = $wb.Worksheets.Item("storage")

$storage.Range("A1").Value = "Y"
$storage.Range("B1").Value = "Z"
$storage.Range("C1").Value = "S"
$storage.Range("D1").Value = "BEL"
$storage.Range("E1").Value = "BEL30"
$storage.Range("F1").Value = "BEL50"
$storage.Range("G1").Value = "BEL70"
Stamp-HeaderStyle $storage.Range("A1:G1")

$storageRows = @(
    @(2050, "BEL_Z", "STOR_L",   "0","0","0","0"),
    @(2050, "BEL_Z", "STOR_M_C", "1000","1000","1000","1000"),
    @(2050, "BEL_Z", "STOR_M_E", "0","0","0","4396.81319461363"),
    @(2050, "BEL_Z", "STOR_S",   "0","0","0","0")
)

$r = 2
foreach ($row in $storageRows) {
    $storage.Cells.Item($r, 1).Value = $row[0]
    $storage.Cells.Item($r, 2).Value = $row[1]
    $storage.Cells.Item($r, 3).Value = $row[2]
    $storage.Cells.Item($r, 4).Value = [double]$row[3]
    $storage.Cells.Item($r, 5).Value = [double]$row[4]
    $storage.Cells.Item($r, 6).Value = [double]$row[5]
    $storage.Cells.Item($r, 7).Value = [double]$row[6]
    Stamp-HeaderStyle $storage.Range("A" + $r + ":C" + $r)
    $r++
}

# ---------------------------------------------------------------------------
# 5) generation (renamed from storage)
# ---------------------------------------------------------------------------
$generation = $wb.Worksheets.Item("generation")

$generation.Range("A1").Value = "Y"
$generation.Range("B1").Value = "G"
$generation.Range("C1").Value = "BEL"
$generation.Range("D1").Value = "BEL30"
$generation.Range("E1").Value = "BEL50"
$generation.Range("F1").Value = "BEL70"
Stamp-HeaderStyle $generation.Range("A1:F1")

$generationRows = @(
    @(2050, "CCGT",    "12877247.3642119","21048421.5158628","46483160.7369532","20133782.3863803"),
    @(2050, "Coal",    "15798960.7915502","11844795.1291264","0","11715276.1679977"),
    @(2050, "Nuclear", "61140089.7944411","37789164.9285595","3610195.26153611","-8.21271789049658e-06"),
    @(2050, "OCGT",    "642645.984971586","910040.464481819","1207569.35621275","1540257.26972145"),
    @(2050, "PV",      "3064318.49854914","3064318.49854947","2950000.35680466","16188373.7107437"),
    @(2050, "WIND_OFF","0","0","0","0"),
    @(2050, "WIND_ON", "8696745.01860184","27578846.5777327","48122010.0627186","55313901.1789302")
)

$r = 2
foreach ($row in $generationRows) {
    $generation.Cells.Item($r, 1).Value = $row[0]
    $generation.Cells.Item($r, 2).Value = $row[1]
    $generation.Cells.Item($r, 3).Value = [double]$row[2]
    $generation.Cells.Item($r, 4).Value = [double]$row[3]
    $generation.Cells.Item($r, 5).Value = [double]$row[4]
    $generation.Cells.Item($r, 6).Value = [double]$row[5]
    Stamp-HeaderStyle $generation.Range("A" + $r + ":B" + $r)
    $r++
}

# ---------------------------------------------------------------------------
# 6) demandshift
# ---------------------------------------------------------------------------
$demandshift.Range("A1").Value = "Z"
$demandshift.Range("B1").Value = "BEL"
$demandshift.Range("C1").Value = "BEL30"
$demandshift.Range("D1").Value = "BEL50"
$demandshift.Range("E1").Value = "BEL70"
Stamp-HeaderStyle $demandshift.Range("A1:E1")

$demandshift.Range("A2").Value = "BEL_Z"
$demandshift.Cells.Item(2, 2).Value = [double]"4107481.76863323"
$demandshift.Cells.Item(2, 3).Value = [double]"5737230.1220309"
$demandshift.Cells.Item(2, 4).Value = [double]"8048022.56185818"
$demandshift.Cells.Item(2, 5).Value = [double]"8747620.05150791"
Stamp-HeaderStyle $demandshift.Range("A2")

# ---------------------------------------------------------------------------
# 7) curtailment (new layout)
# ---------------------------------------------------------------------------
$curtailment.Range("A1").Value = "Y"
$curtailment.Range("B1").Value = "curt"
$curtailment.Range("C1").Value = "curt30"
$curtailment.Range("D1").Value = "curt50"
$curtailment.Range("E1").Value = "curt70"
Stamp-HeaderStyle $curtailment.Range("A1:E1")

$curtailment.Cells.Item(2, 1).Value = 2050
$curtailment.Cells.Item(2, 2).Value = [double]"0"
$curtailment.Cells.Item(2, 3).Value = [double]"9692.70179831109"
$curtailment.Cells.Item(2, 4).Value = [double]"1928414.36552252"
$curtailment.Cells.Item(2, 5).Value = [double]"4751047.75847795"
Stamp-HeaderStyle $curtailment.Range("A2")

# ---------------------------------------------------------------------------
# 8) storc
# ---------------------------------------------------------------------------
$storc.Range("A1").Value = "Y"
$storc.Range("B1").Value = "BEL"
$storc.Range("C1").Value = "BEL30"
$storc.Range("D1").Value = "BEL50"
$storc.Range("E1").Value = "BEL70"
Stamp-HeaderStyle $storc.Range("A1:E1")

$storc.Cells.Item(2, 1).Value = 2050
$storc.Cells.Item(2, 2).Value = [double]"300906.110457436"
$storc.Cells.Item(2, 3).Value = [double]"366814.106616291"
$storc.Cells.Item(2, 4).Value = [double]"915659.740745972"
$storc.Cells.Item(2, 5).Value = [double]"10981934.914065"
Stamp-HeaderStyle $storc.Range("A2")

# ---------------------------------------------------------------------------
# 9) stordisc
# ---------------------------------------------------------------------------
$stordisc.Range("A1").Value = "Y"
$stordisc.Range("B1").Value = "BEL"
$stordisc.Range("C1").Value = "BEL30"
$stordisc.Range("D1").Value = "BEL50"
$stordisc.Range("E1").Value = "BEL70"
Stamp-HeaderStyle $stordisc.Range("A1:E1")

$stordisc.Cells.Item(2, 1).Value = 2050
$stordisc.Cells.Item(2, 2).Value = [double]"225679.582843074"
$stordisc.Cells.Item(2, 3).Value = [double]"275110.579962224"
$stordisc.Cells.Item(2, 4).Value = [double]"686744.805559474"
$stordisc.Cells.Item(2, 5).Value = [double]"8236451.18555302"
Stamp-HeaderStyle $stordisc.Range("A2")

# ---------------------------------------------------------------------------
# 10) maxshift
# ---------------------------------------------------------------------------
$maxshift.Range("A1").Value = "P"
$maxshift.Range("B1").Value = "BEL"
$maxshift.Range("C1").Value = "BEL30"
$maxshift.Range("D1").Value = "BEL50"
$maxshift.Range("E1").Value = "BEL70"
Stamp-HeaderStyle $maxshift.Range("A1:E1")

$maxshiftRows = @(
    @(1,"2129.87465686798","3998.4748798743","4889.57276278364","4889.57276278365"),
    @(2,"2265.99565424034","2312.8959089581","2321.30935452233","2321.30935452233"),
    @(3,"2247.36667348608","2312.8959089581","2321.30935452233","2321.30935452233"),
    @(4,"1373.524963539","1454.59514560029","1454.5951456003","1454.5951456003"),
    @(5,"1382.56860452993","1454.5951456003","1454.5951456003","1454.5951456003"),
    @(6,"3521.80391950026","3549.49166202945","3596.11976307164","3596.11976307164"),
    @(7,"3335.69743176395","3510.51917214959","3610.51917214959","3610.51917214959"),
    @(8,"2903.22867105034","4857.62891670642","4857.62891670642","4857.62891670642")
)

$r = 2
foreach ($row in $maxshiftRows) {
    $maxshift.Cells.Item($r, 1).Value = $row[0]
    $maxshift.Cells.Item($r, 2).Value = [double]$row[1]
    $maxshift.Cells.Item($r, 3).Value = [double]$row[2]
    $maxshift.Cells.Item($r, 4).Value = [double]$row[3]
    $maxshift.Cells.Item($r, 5).Value = [double]$row[4]
    Stamp-HeaderStyle $maxshift.Range("A" + $r)
    $r++
}

# ---------------------------------------------------------------------------
# 11) restore the original active sheet
# ---------------------------------------------------------------------------
$reserves.Activate()
